$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = "ba"
$ws.Range("J3").Value = "Appreciation"
$ws.Range("I7").Value = "sv"
$ws.Range("J7").Value = "Statement-opinion"
$ws.Range("I16").Value = "%"
$ws.Range("J16").Value = "Uninterpretable"
$ws.Range("I19").Value = "sv"
$ws.Range("J19").Value = "Statement-opinion"
$ws.Range("I20").Value = "sd"
$ws.Range("J20").Value = "Statement-non-opinion"
$ws.Range("I23").Value = "ba"
$ws.Range("J23").Value = "Appreciation"
$ws.Range("I36").Value = "sv"
$ws.Range("J36").Value = "Statement-opinion"
$ws.Range("I38").Value = "sv"
$ws.Range("J38").Value = "Statement-opinion"
$ws.Range("I54").Value = "aa"
$ws.Range("J54").Value = "Agree/Accept"
$ws.Range("I55").Value = "sd"
$ws.Range("J55").Value = "Statement-non-opinion"
$ws.Range("I57").Value = "ba"
$ws.Range("J57").Value = "Appreciation"
$ws.Range("I61").Value = "aa"
$ws.Range("J61").Value = "Agree/Accept"
$ws.Range("I62").Value = "aa"
$ws.Range("J62").Value = "Agree/Accept"
$ws.Range("I70").Value = "sv"
$ws.Range("J70").Value = "Statement-opinion"
$ws.Range("I74").Value = "sd"
$ws.Range("J74").Value = "Statement-non-opinion"
$ws.Range("I75").Value = "sd"
$ws.Range("J75").Value = "Statement-non-opinion"
$ws.Range("I78").Value = "aa"
$ws.Range("J78").Value = "Agree/Accept"
$ws.Range("I88").Value = "ba"
$ws.Range("J88").Value = "Appreciation"
$ws.Range("I92").Value = "aa"
$ws.Range("J92").Value = "Agree/Accept"
$ws.Range("I95").Value = "sd"
$ws.Range("J95").Value = "Statement-non-opinion"
$ws.Range("I98").Value = "sd"
$ws.Range("J98").Value = "Statement-non-opinion"
$ws.Range("I101").Value = "sd"
$ws.Range("J101").Value = "Statement-non-opinion"
$ws.Range("I105").Value = "sd"
$ws.Range("J105").Value = "Statement-non-opinion"
$ws.Range("I109").Value = "sv"
$ws.Range("J109").Value = "Statement-opinion"
$ws.Range("I113").Value = "sd"
$ws.Range("J113").Value = "Statement-non-opinion"
$ws.Range("I116").Value = "sd"
$ws.Range("J116").Value = "Statement-non-opinion"
$ws.Range("I118").Value = "%"
$ws.Range("J118").Value = "Uninterpretable"
$ws.Range("I133").Value = "sv"
$ws.Range("J133").Value = "Statement-opinion"
$ws.Range("I135").Value = "sd"
$ws.Range("J135").Value = "Statement-non-opinion"
$ws.Range("I137").Value = "sd"
$ws.Range("J137").Value = "Statement-non-opinion"
$ws.Range("I139").Value = "sv"
$ws.Range("J139").Value = "Statement-opinion"
$ws.Range("I148").Value = "qy"
$ws.Range("J148").Value = "Yes-No-Question"
$ws.Range("I159").Value = "sd"
$ws.Range("J159").Value = "Statement-non-opinion"
$ws.Range("I177").Value = "sv"
$ws.Range("J177").Value = "Statement-opinion"
$ws.Range("I178").Value = "%"
$ws.Range("J178").Value = "Uninterpretable"
$ws.Range("I182").Value = "aa"
$ws.Range("J182").Value = "Agree/Accept"
$ws.Range("I183").Value = "aa"
$ws.Range("J183").Value = "Agree/Accept"
$ws.Range("I198").Value = "sd"
$ws.Range("J198").Value = "Statement-non-opinion"
$ws.Range("I206").Value = "qy"
$ws.Range("J206").Value = "Yes-No-Question"
$ws.Range("I221").Value = "aa"
$ws.Range("J221").Value = "Agree/Accept"
$ws.Range("I222").Value = "sd"
$ws.Range("J222").Value = "Statement-non-opinion"
$ws.Range("I231").Value = "sv"
$ws.Range("J231").Value = "Statement-opinion"
$ws.Range("I234").Value = "ba"
$ws.Range("J234").Value = "Appreciation"
$ws.Range("I235").Value = "sd"
$ws.Range("J235").Value = "Statement-non-opinion"
$ws.Range("I238").Value = "aa"
$ws.Range("J238").Value = "Agree/Accept"
$ws.Range("I245").Value = "b"
$ws.Range("J245").Value = "Acknowledge (Backchannel)"
$ws.Range("I264").Value = "sd"
$ws.Range("J264").Value = "Statement-non-opinion"
$ws.Range("I266").Value = "sd"
$ws.Range("J266").Value = "Statement-non-opinion"
$ws.Range("I273").Value = "aa"
$ws.Range("J273").Value = "Agree/Accept"
$ws.Range("I274").Value = "b"
$ws.Range("J274").Value = "Acknowledge (Backchannel)"
$ws.Range("I282").Value = "sd"
$ws.Range("J282").Value = "Statement-non-opinion"
$ws.Range("I295").Value = "sd"
$ws.Range("J295").Value = "Statement-non-opinion"
$ws.Range("I318").Value = "sv"
$ws.Range("J318").Value = "Statement-opinion"
$ws.Range("I320").Value = "sv"
$ws.Range("J320").Value = "Statement-opinion"
